$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new requirement row (row 19) below the existing data (row 18),
# reusing the same text as row 18 ("Before a deployment is made..." / No /
# N/A / N/A) with the next sequence number (18).
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = $ws.Range("B18").Text
$ws.Range("C19").Value = $ws.Range("C18").Text
$ws.Range("D19").Value = $ws.Range("D18").Text
$ws.Range("E19").Value = $ws.Range("E18").Text
$ws.Range("F19").Value = $ws.Range("F18").Text
$ws.Range("G19").Value = $ws.Range("G18").Text
$ws.Range("H19").Value = $ws.Range("H18").Text

# Match the formatting used by the rest of the data rows (wrapped text,
# top-aligned, boxed in a thin border) and the row height used for row 18.
$ws.Range("A19:H19").VerticalAlignment = -4160
$ws.Range("A19:H19").WrapText = $true
$ws.Range("A19:H19").Borders.LineStyle = 1
$ws.Rows.Item(19).RowHeight = $ws.Rows.Item(18).RowHeight

# Leave the new row's last cell selected, as it was when editing finished.
$ws.Range("B19").Select()
